$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. "1.00", "6.10")
# are not auto-converted to numbers, matching the inline-string cells in the source file.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.251.88"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.584.78"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "570.10"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "143.50"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.594.40"
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  +9.60%  "
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "3.042.58"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "59.305.08"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "22.52"
$ws.Range("E16").Value = "  +7.93%  "
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "2.588.82"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "336.32"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "64.21"
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("E25").Value = "  +5.97%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "7.26"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "6.10"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "156.72"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "0.889"
$ws.Range("E36").Value = "  +8.03%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.881"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.14"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "36.83"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "295.65"
$ws.Range("E41").Value = "  +4.53%  "
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "0.0976"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "19.22"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("D48").Value = "10.62"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "124.57"
$ws.Range("E49").Value = "  +4.56%  "
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E51").Value = "  +4.16%  "
